$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Status columns (zh-cn / de-de) move from "Handed back: in sync with en-US" to "Ready for handoff"
# and the "Latest HO Xliff Generate Date" timestamp is refreshed.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 23:01:04"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ---
# Status column moves to "Ready for handoff" and the handoff datetime is refreshed.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 23:00:56"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ---
# Status column moves to "Ready for handoff" and the handoff datetime is refreshed.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 23:01:04"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
